# Update the "Example-AttachFile" worksheet so that the MultiFormParams
# value cell (J2) contains the additional multipart form information
# (serverUrls instead of serverurl, plus a reportTitle parameter).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example-AttachFile")

$ws.Range("J2").Value = "filestream=sample.json;serverUrls=https://live.virtualandemo.com/api;dataload=APITEST.json;execute=true;type=VIRTUALAN;reportTitle=DemoTestReport"

# Move the active selection to J3, matching where the cursor ended up
# after the edit.
$ws.Activate()
$ws.Range("J3").Select()
